# Comprobación de los posibles errores de lectura
# Replace the sample data (Raul/Claudia/Laura/Carmen) with new test rows
# (Gggg/Hhhh/Iiii/Jjjj), including NIF values and emails, and move the
# active cell selection from D5 to C5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 2: Gggg / 88888888G / gg@uniovi.es
$ws.Range("A2").Value = "Gggg"
$ws.Range("B2").Value = "88888888G"
$ws.Range("C2").Value = "gg@uniovi.es"

# Row 3: Hhhh / 999999999 (numeric) / hh@uniovi.es
$ws.Range("A3").Value = "Hhhh"
$ws.Range("B3").Value = 999999999
$ws.Range("C3").Value = "hh@uniovi.es"

# Row 4: Iiii / 131313131I / ii@uniovi.es
$ws.Range("A4").Value = "Iiii"
$ws.Range("B4").Value = "131313131I"
$ws.Range("C4").Value = "ii@uniovi.es"

# Row 5: Jjjj / 14141414J / jj@uniovi.es
$ws.Range("A5").Value = "Jjjj"
$ws.Range("B5").Value = "14141414J"
$ws.Range("C5").Value = "jj@uniovi.es"

# Move the selected/active cell from D5 to C5
$ws.Range("C5").Select()
